# Apply the syllabus update described in the commit:
#  - bump the term from Fall 2016 to Fall 2018
#  - move the lecture location from Barrows Hall 122 to Campbell Hall 501
#  - tweak the oral-final bullet's wording ("and" -> ",")
#  - note that a jupyter installation is an optional material, ahead of LaTeX

$d = $word.ActiveDocument

# 1) Fall 2016 -> Fall 2018
$d.Content.Find.Execute("Fall 2016", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Fall 2018", 2) | Out-Null

# 2) Barrows Hall 122 -> Campbell Hall 501
$d.Content.Find.Execute("Barrows Hall 122", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Campbell Hall 501", 2) | Out-Null

# 3) "... prelims/quals and giving scientific talks" -> "... prelims/quals , giving scientific talks"
$d.Content.Find.Execute("quals and giving scientific talks", $true, $false, $false, $false, $false, `
    $true, 1, $false, "quals , giving scientific talks", 2) | Out-Null

# 4) Materials bullet: "a LaTeX installation" -> "optional: a jupyter installation and a LaTeX installation"
$d.Content.Find.Execute("a LaTeX installation", $true, $false, $false, $false, $false, `
    $true, 1, $false, "optional: a jupyter installation and a LaTeX installation", 2) | Out-Null
